$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B ("ASIN") to host the new
# "Week_Start_Date" column. Everything from the old column B onward
# (ASIN, MyForecast, Amazon Mean/P70/P80/P90 Forecast, Product Title,
# is_holiday_week) shifts one column to the right.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week labels in column A: drop the leading zero (W01 -> W1, ... W16 stays W16).
$weekLabels = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")

# Monday week-start dates for column B, one per data row (rows 2-17).
$weekStartDates = @(
    "2025-01-05","2025-01-12","2025-01-19","2025-01-26",
    "2025-02-02","2025-02-09","2025-02-16","2025-02-23",
    "2025-03-02","2025-03-09","2025-03-16","2025-03-23",
    "2025-03-30","2025-04-06","2025-04-13","2025-04-20"
)

for ($i = 0; $i -lt $weekLabels.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $weekLabels[$i]
    # Force the date-looking text to stay a literal string (not auto-parsed
    # into a date serial number) by entering it the way a user would -
    # with a leading apostrophe - then drop the resulting "text quote
    # prefix" formatting so the cell keeps the plain default style.
    $ws.Cells.Item($r, 2).Value = "'" + $weekStartDates[$i]
    $ws.Cells.Item($r, 2).ClearFormats()
}

$wb.Save()
